# "few changes in scenario LP"
#
# - Rename sheet 1 (currently "Properties")  -> "Properties_more"
# - Rename sheet 2 (currently "Properties2") -> "Properties"
# - Make the (renamed) "Properties" sheet the active/selected tab
# - Move that sheet's selection from E40 to F60

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "Properties_more"
$ws2.Name = "Properties"

$ws2.Activate()
$ws2.Range("F60").Select()
